$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 17 de Agosto de 2020 a las 22:57"

# Refresh country case counts. A handful of rows also swap position
# with their neighbour because the sheet stays sorted by column B
# (Casos totales) descending.

$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 5595414
$ws.Range("C4").Value = 29360
$ws.Range("D4").Value = 2948201
$ws.Range("E4").Value = 2473736
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 350
$ws.Range("H4").Value = 173477

$ws.Range("A6").Value = "India"
$ws.Range("B6").Value = 2701604
$ws.Range("C6").Value = 54288
$ws.Range("D6").Value = 1976248
$ws.Range("E6").Value = 673431
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 880
$ws.Range("H6").Value = 51925

$ws.Range("A8").Value = "Sudafrica"
$ws.Range("B8").Value = 589886
$ws.Range("C8").Value = 2541
$ws.Range("D8").Value = 477671
$ws.Range("E8").Value = 100233
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 143
$ws.Range("H8").Value = 11982

$ws.Range("A22").Value = "Alemania"
$ws.Range("B22").Value = 226537
$ws.Range("C22").Value = 1540
$ws.Range("D22").Value = 202900
$ws.Range("E22").Value = 14341
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 6
$ws.Range("H22").Value = 9296

$ws.Range("A33").Value = "Israel"
$ws.Range("B33").Value = 94751
$ws.Range("C33").Value = 2071
$ws.Range("D33").Value = 70291
$ws.Range("E33").Value = 23768
$ws.Range("F33").Value = 0
$ws.Range("G33").Value = 7
$ws.Range("H33").Value = 692

$ws.Range("A68").Value = "Costa Rica"
$ws.Range("B68").Value = 29084
$ws.Range("C68").Value = 619
$ws.Range("D68").Value = 9233
$ws.Range("E68").Value = 19547
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 10
$ws.Range("H68").Value = 304

$ws.Range("A76").Value = "Costa de Marfil"
$ws.Range("B76").Value = 17107
$ws.Range("C76").Value = 81
$ws.Range("D76").Value = 13990
$ws.Range("E76").Value = 3007
$ws.Range("F76").Value = 0
$ws.Range("G76").Value = 0
$ws.Range("H76").Value = 110

$ws.Range("A93").Value = "Guinea"
$ws.Range("B93").Value = 8620
$ws.Range("C93").Value = 138
$ws.Range("D93").Value = 7472
$ws.Range("E93").Value = 1097
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = 51

$ws.Range("A94").Value = "Gabon"
$ws.Range("B94").Value = 8270
$ws.Range("C94").Value = 45
$ws.Range("D94").Value = 6404
$ws.Range("E94").Value = 1813
$ws.Range("F94").Value = 0
$ws.Range("G94").Value = 2
$ws.Range("H94").Value = 53

$ws.Range("A97").Value = "Haiti"
$ws.Range("B97").Value = 7897
$ws.Range("C97").Value = 18
$ws.Range("D97").Value = 5235
$ws.Range("E97").Value = 2466
$ws.Range("F97").Value = 0
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 196

$ws.Range("A100").Value = "Luxemburgo"
$ws.Range("B100").Value = 7469
$ws.Range("C100").Value = 11
$ws.Range("D100").Value = 6739
$ws.Range("E100").Value = 606
$ws.Range("F100").Value = 0
$ws.Range("G100").Value = 1
$ws.Range("H100").Value = 124

$ws.Range("A104").Value = "Maldivas"
$ws.Range("B104").Value = 5909
$ws.Range("C104").Value = 124
$ws.Range("D104").Value = 3488
$ws.Range("E104").Value = 2398
$ws.Range("F104").Value = 0
$ws.Range("G104").Value = 1
$ws.Range("H104").Value = 23

$ws.Range("A105").Value = "Republica de Yibuti"
$ws.Range("B105").Value = 5372
$ws.Range("C105").Value = 3
$ws.Range("D105").Value = 5211
$ws.Range("E105").Value = 102
$ws.Range("F105").Value = 0
$ws.Range("G105").Value = 0
$ws.Range("H105").Value = 59

$ws.Range("A107").Value = "Malaui"
$ws.Range("B107").Value = 5125
$ws.Range("C107").Value = 53
$ws.Range("D107").Value = 2690
$ws.Range("E107").Value = 2273
$ws.Range("F107").Value = 0
$ws.Range("G107").Value = 1
$ws.Range("H107").Value = 162

$ws.Range("A110").Value = "Republica de Africa Central"
$ws.Range("B110").Value = 4667
$ws.Range("C110").Value = 15
$ws.Range("D110").Value = 1748
$ws.Range("E110").Value = 2858
$ws.Range("F110").Value = 0
$ws.Range("G110").Value = 0
$ws.Range("H110").Value = 61

$ws.Range("A115").Value = "Suazilandia"
$ws.Range("B115").Value = 3894
$ws.Range("C115").Value = 55
$ws.Range("D115").Value = 2371
$ws.Range("E115").Value = 1450
$ws.Range("F115").Value = 0
$ws.Range("G115").Value = 3
$ws.Range("H115").Value = 73

$ws.Range("A120").Value = "Cabo Verde"
$ws.Range("B120").Value = 3203
$ws.Range("C120").Value = 24
$ws.Range("D120").Value = 2336
$ws.Range("E120").Value = 831
$ws.Range("F120").Value = 0
$ws.Range("G120").Value = 1
$ws.Range("H120").Value = 36

$ws.Range("A121").Value = "Mayotte"
$ws.Range("B121").Value = 3160
$ws.Range("C121").Value = 41
$ws.Range("D121").Value = 2964
$ws.Range("E121").Value = 157
$ws.Range("F121").Value = 0
$ws.Range("G121").Value = 0
$ws.Range("H121").Value = 39

$ws.Range("A127").Value = "Ruanda"
$ws.Range("B127").Value = 2540
$ws.Range("C127").Value = 87
$ws.Range("D127").Value = 1661
$ws.Range("E127").Value = 871
$ws.Range("F127").Value = 0
$ws.Range("G127").Value = 0
$ws.Range("H127").Value = 8

$ws.Range("A128").Value = "Sudan del Sur"
$ws.Range("B128").Value = 2490
$ws.Range("C128").Value = 1
$ws.Range("D128").Value = 1175
$ws.Range("E128").Value = 1268
$ws.Range("F128").Value = 0
$ws.Range("G128").Value = 0
$ws.Range("H128").Value = 47

$ws.Range("A132").Value = "Tunez"
$ws.Range("B132").Value = 2185
$ws.Range("C132").Value = 78
$ws.Range("D132").Value = 1362
$ws.Range("E132").Value = 767
$ws.Range("F132").Value = 0
$ws.Range("G132").Value = 2
$ws.Range("H132").Value = 56

$ws.Range("A133").Value = "Guinea-Bisau"
$ws.Range("B133").Value = 2117
$ws.Range("C133").Value = 0
$ws.Range("D133").Value = 1015
$ws.Range("E133").Value = 1069
$ws.Range("F133").Value = 0
$ws.Range("G133").Value = 0
$ws.Range("H133").Value = 33

$ws.Range("A137").Value = "Angola"
$ws.Range("B137").Value = 1935
$ws.Range("C137").Value = 29
$ws.Range("D137").Value = 632
$ws.Range("E137").Value = 1215
$ws.Range("F137").Value = 0
$ws.Range("G137").Value = 0
$ws.Range("H137").Value = 88

$ws.Range("A138").Value = "Yemen"
$ws.Range("B138").Value = 1882
$ws.Range("C138").Value = 13
$ws.Range("D138").Value = 1045
$ws.Range("E138").Value = 302
$ws.Range("F138").Value = 0
$ws.Range("G138").Value = 5
$ws.Range("H138").Value = 535

$ws.Range("A139").Value = "Gambia"
$ws.Range("B139").Value = 1872
$ws.Range("C139").Value = 0
$ws.Range("D139").Value = 401
$ws.Range("E139").Value = 1408
$ws.Range("F139").Value = 0
$ws.Range("G139").Value = 0
$ws.Range("H139").Value = 63

$ws.Range("A151").Value = "Burkina Faso"
$ws.Range("B151").Value = 1280
$ws.Range("C151").Value = 13
$ws.Range("D151").Value = 1018
$ws.Range("E151").Value = 207
$ws.Range("F151").Value = 0
$ws.Range("G151").Value = 0
$ws.Range("H151").Value = 55

$ws.Range("A152").Value = "Liberia"
$ws.Range("B152").Value = 1277
$ws.Range("C152").Value = 20
$ws.Range("D152").Value = 803
$ws.Range("E152").Value = 392
$ws.Range("F152").Value = 0
$ws.Range("G152").Value = 0
$ws.Range("H152").Value = 82

$ws.Range("A154").Value = "Togo"
$ws.Range("B154").Value = 1154
$ws.Range("C154").Value = 7
$ws.Range("D154").Value = 858
$ws.Range("E154").Value = 269
$ws.Range("F154").Value = 0
$ws.Range("G154").Value = 0
$ws.Range("H154").Value = 27

$ws.Range("A155").Value = "Aruba"
$ws.Range("B155").Value = 1121
$ws.Range("C155").Value = 19
$ws.Range("D155").Value = 212
$ws.Range("E155").Value = 905
$ws.Range("F155").Value = 0
$ws.Range("G155").Value = 0
$ws.Range("H155").Value = 4

$ws.Range("A156").Value = "Jamaica"
$ws.Range("B156").Value = 1113
$ws.Range("C156").Value = 7
$ws.Range("D156").Value = 764
$ws.Range("E156").Value = 335
$ws.Range("F156").Value = 0
$ws.Range("G156").Value = 0
$ws.Range("H156").Value = 14

$ws.Range("A199").Value = "Curazao"
$ws.Range("B199").Value = 36
$ws.Range("C199").Value = 2
$ws.Range("D199").Value = 31
$ws.Range("E199").Value = 4
$ws.Range("F199").Value = 0
$ws.Range("G199").Value = 0
$ws.Range("H199").Value = 1

$ws.Range("A213").Value = "Montserrat"
$ws.Range("B213").Value = 13
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 12
$ws.Range("E213").Value = 0
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 1

$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 13
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 0
